$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "In Translation"
#    (shared-string cell, referenced from the Overview / zh-cn / de-de sheets)
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$overview.Range("E4").Value = "In Translation"
$overview.Range("F4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"
$dede.Range("C4").Value = "In Translation"

# ---------------------------------------------------------------------------
# 2. Narrow the "zh-cn" / "de-de" (Overview) and "Status" (zh-cn / de-de)
#    columns to their post-refresh report width.
# ---------------------------------------------------------------------------
$overview.Range("E1:F1").EntireColumn.ColumnWidth = 12.5
$zhcn.Range("C1").EntireColumn.ColumnWidth = 12.5
$dede.Range("C1").EntireColumn.ColumnWidth = 12.5
